$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (F1, G1, H1)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy the style of the existing header (A1) onto the new headers so they
# match the bold / bordered / centered look of the rest of row 1.
$ws.Range("A1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)

# Boolean outlier-flag values for columns F (KNN), G (SVM) and H (RF),
# one row per data row (rows 2-25).
$outlierData = @(
    @($false, $false, $false),  # row 2
    @($false, $false, $false),  # row 3
    @($false, $false, $false),  # row 4
    @($false, $false, $false),  # row 5
    @($false, $false, $false),  # row 6
    @($false, $false, $false),  # row 7
    @($false, $false, $false),  # row 8
    @($true,  $true,  $true ),  # row 9
    @($false, $false, $false),  # row 10
    @($false, $false, $false),  # row 11
    @($false, $false, $false),  # row 12
    @($false, $false, $false),  # row 13
    @($false, $false, $false),  # row 14
    @($false, $false, $false),  # row 15
    @($false, $false, $false),  # row 16
    @($false, $false, $false),  # row 17
    @($false, $true,  $true ),  # row 18
    @($false, $false, $false),  # row 19
    @($false, $false, $false),  # row 20
    @($false, $false, $false),  # row 21
    @($false, $false, $false),  # row 22
    @($false, $false, $false),  # row 23
    @($true,  $true,  $true ),  # row 24
    @($true,  $false, $true )   # row 25
)

for ($i = 0; $i -lt $outlierData.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $outlierData[$i][0]
    $ws.Cells.Item($row, 7).Value = $outlierData[$i][1]
    $ws.Cells.Item($row, 8).Value = $outlierData[$i][2]
}
